# Applies the changes described by the commit:
#  - "row with a boolean" -> "row with a boolean formula" (Sheet1!B7)
#  - "row with a formula" -> "row with formulas" (Sheet1!B8)
#  - Sheet1!C8 AVERAGE formula drops the 128 argument (new result 21)
#  - Sheet1!D8 gains a new CONCATENATE formula ("string cat")
#  - Sheet3 renamed to "2 of 3" and filled with a small data table, and
#    made the active sheet
#  - selection/active-cell bookkeeping on Sheet1 / Sheet3
#  - minor row-height bookkeeping (customHeight) on existing rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B7").Value = "row with a boolean formula"
$ws1.Range("B8").Value = "row with formulas"

$ws1.Range("C8").Formula = "=AVERAGE(2,4,8,16,32,64)"
$ws1.Range("D8").Formula = "=CONCATENATE(""string "",""cat"")"
$ws1.Range("B8:D8").Font.Name = "Arial"

foreach ($r in @(1,2,3,5,6,7,8,9,10)) {
    $ws1.Rows.Item($r).RowHeight = 12.8
}

$ws1.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------
# 2ndsheet - no content changes, just row-height bookkeeping
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2ndsheet")
foreach ($r in @(1,2,3)) {
    $ws2.Rows.Item($r).RowHeight = 12.8
}

# ---------------------------------------------------------------------
# Sheet3 -> "2 of 3"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "2 of 3"

$ws3.Range("A1").Value = "a sheet"
$ws3.Range("B1").Value = "with"
$ws3.Range("C1").Value = 4
$ws3.Range("D1").Value = "columns"
$ws3.Range("A1:D1").Font.Bold = $true

$ws3.Range("A2").Value = "row 1,col 1"
$ws3.Range("B2").Value = "row 1,col 2"
$ws3.Range("C2").Value = "row 1,col 3"
$ws3.Range("D2").Value = "row 1,col 4"

$ws3.Range("A3").Value = "row 2,col 1"
$ws3.Range("B3").Value = "row 2,col 2"
$ws3.Range("C3").Value = "row 2,col 3"
$ws3.Range("D3").Value = "row 2,col 4"

$ws3.Range("A2:D3").Font.Name = "Arial"

# Make Sheet3 ("2 of 3") the active sheet with the right selection.
$ws3.Activate() | Out-Null
$ws3.Range("C16").Select() | Out-Null
